# Update countries & provincias Spain
# Applies the COVID data refresh captured in the commit diff:
#  - Swap display order of "Santa Lucia" / "Timor Oriental" rows
#  - Bump the "Datos actualizados ..." timestamp from 14:22 to 15:39
#  - Refresh the numeric stats (Casos totales, Nuevos casos, Casos activos,
#    Casos criticos, Muertes) for the affected country rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string reorder: Timor Oriental now precedes Santa Lucia ---
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# --- Timestamp refresh ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 15:39"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 5842153
$ws.Range("C4").Value = 725
$ws.Range("E4").Value = 2513882

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 3583308
$ws.Range("C5").Value = 610
$ws.Range("E5").Value = 759383
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 114287

# --- Row 14: Iran ---
$ws.Range("B14").Value = 358905
$ws.Range("C14").Value = 2113
$ws.Range("D14").Value = 309464
$ws.Range("E14").Value = 28798
$ws.Range("G14").Value = 141
$ws.Range("H14").Value = 20643

# --- Row 17: Arabia Saudita ---
$ws.Range("B17").Value = 307479
$ws.Range("C17").Value = 1109
$ws.Range("D17").Value = 280143
$ws.Range("E17").Value = 23687
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 3649

# --- Row 43: Bielorrusia ---
$ws.Range("B43").Value = 70468
$ws.Range("C43").Value = 183
$ws.Range("D43").Value = 68839
$ws.Range("E43").Value = 987
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 642

# --- Row 46: Paises Bajos ---
$ws.Range("B46").Value = 66554
$ws.Range("C46").Value = 457

# --- Row 62: Uzbekistan ---
$ws.Range("B62").Value = 38870
$ws.Range("C62").Value = 458
$ws.Range("D62").Value = 34894
$ws.Range("E62").Value = 3705
$ws.Range("G62").Value = 6
$ws.Range("H62").Value = 271

# --- Row 64: Azerbaiyan ---
$ws.Range("B64").Value = 35274
$ws.Range("C64").Value = 169
$ws.Range("D64").Value = 32993
$ws.Range("E64").Value = 1763
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 518

# --- Row 69: Serbia ---
$ws.Range("B69").Value = 30657
$ws.Range("C69").Value = 109
$ws.Range("E69").Value = 1684
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 698

# --- Row 83: Republica de Macedonia ---
$ws.Range("B83").Value = 13595
$ws.Range("C83").Value = 137
$ws.Range("D83").Value = 10110
$ws.Range("E83").Value = 2921
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 564

# --- Row 125: Sri Lanka ---
$ws.Range("B125").Value = 2951
$ws.Range("C125").Value = 4
$ws.Range("E125").Value = 134

# --- Row 144: Malta ---
$ws.Range("B144").Value = 1612
$ws.Range("C144").Value = 35
$ws.Range("D144").Value = 934
$ws.Range("E144").Value = 668
